$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Venta 2"
$ws.Range("B3").Value = 200.0
